$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# --- Make room in the rows layout ---
# Current layout (before):
#   45 pythonPath | C:\Python312
#   46 pythonLibraryPath | C:\Python312\python312.dll
#   47 HashScriptPath | ...hash_policyholder_id.py   (ht=30)
#   48 (blank)
#   49 DbPath | ...claims_encrypted.db                (ht=30)
#   50 SQLiteDSN | ClaimsAutomation
#
# Target layout (after):
#   44 Subject_DatabaseError | Database Connection Error!
#   45 Body_DatabaseError | <long body>                (ht=225)
#   46 (blank)
#   47 pythonPath | C:\Python312
#   48 pythonLibraryPath | C:\Python312\python312.dll
#   49 ValidateScriptPath | ...validate_claims.py       (ht=30)
#   50 HashScriptPath | ...hash_policyholder_id.py      (ht=30)
#   51 (blank)
#   52 DbPath | ...claims_encrypted.db                  (ht=30)
#   53 SQLiteDSN | ClaimsAutomation

# Insert a row above HashScriptPath (row 47) to make space for ValidateScriptPath.
$ws.Rows.Item(47).Insert()

# Insert two rows above pythonPath (row 45) to make space for the new
# Subject_DatabaseError / Body_DatabaseError pair.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# --- Fill in the new content ---
# (ValidateScriptPath is written first so the shared-string table allocation
# order matches the authored workbook.)

# ValidateScriptPath
$ws.Range("A49").Value = "ValidateScriptPath"
$ws.Range("B49").Value = "C:\Users\MorisMwaiWachira\Desktop\MorisMwai_RPA_Assignment\Python_Scripts\validate_claims.py"
$ws.Range("B49").WrapText = $true
$ws.Rows.Item(49).RowHeight = 30

# Subject_DatabaseError / Database Connection Error!
$ws.Range("A44").Value = "Subject_DatabaseError"
$ws.Range("B44").Value = "Database Connection Error!"
$ws.Range("B44").WrapText = $true

# Body_DatabaseError / long body text
$ws.Range("A45").Value = "Body_DatabaseError"
$ws.Range("B45").Value = "Hello,`nAn error occurred while connecting to or interacting with the database. Please review the details below:`nException Source: @Source  `nException Message: @Message  `nThis might be due to incorrect DSN configuration, network issues, or driver-specific errors.`nA screenshot of the error has been attached for reference.`nThank you,  `nRobot :)"
$ws.Range("B45").WrapText = $true
$ws.Rows.Item(45).RowHeight = 225

# --- View state: make Constants the active/selected sheet, matching the diff ---
$ws.Activate()
$ws.Range("B49").Select()
